# Actualizacion automatica del mapa (2025-10-29 11:16:32)
# Adds two new rows (87 and 88) of data to the "Optical_Power" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT even though it looks
# like a number or a date (Excel would otherwise auto-convert it).
function Set-TextValue {
    param(
        [int]$Row,
        [int]$Col,
        [string]$Value
    )
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    # Reset the style back to the workbook default so no extra formatting
    # (besides the "text" type) is introduced.
    $cell.Style = "Normal"
}

# Helper: write a plain text value (safe - Excel will not reinterpret it).
function Set-StringValue {
    param(
        [int]$Row,
        [int]$Col,
        [string]$Value
    )
    $ws.Cells.Item($Row, $Col).Value = $Value
}

# Helper: write a numeric value.
function Set-NumberValue {
    param(
        [int]$Row,
        [int]$Col,
        [double]$Value
    )
    $ws.Cells.Item($Row, $Col).Value = $Value
}

$rows = @(
    @{
        Row = 87
        Caso = "7665"
        Fecha = "10/28/2025"
        Direccion = "ARAOZ 2313"
        Comuna = "14"
        OT = "810461115"
        Proveedor = "Optical Power"
        Estado = "Pendiente"
        Observaciones = "Picada y cable cortado"
        Attachments = 1
        Tipo = "Cambio"
        Equipo = "Sin equipos"
        TipoElemento = "Pasante"
        CoordX = -58.417634
        CoordY = -34.587439
        Operacion = "Palermo"
        Zona = "Capital Sur"
        PD = "VCR-G"
        N2 = "Fuera de Poligono OVL"
    },
    @{
        Row = 88
        Caso = "7619"
        Fecha = "10/29/2025"
        Direccion = "QUITO 4180"
        Comuna = "5"
        OT = "810471618"
        Proveedor = "Optical Power"
        Estado = "Pendiente"
        Observaciones = "Picada"
        Attachments = 1
        Tipo = "Cambio"
        Equipo = "Sin equipos"
        TipoElemento = "Pasante"
        CoordX = -58.425596
        CoordY = -34.617038
        Operacion = "Almagro"
        Zona = "Capital Sur"
        PD = "ALM-C"
        N2 = "Fuera de Poligono OVL"
    }
)

foreach ($r in $rows) {
    $row = $r.Row

    # A: Caso (numeric-looking text)
    Set-TextValue $row 1 $r.Caso
    # B: F. De Reclamo (date stored as text)
    Set-TextValue $row 2 $r.Fecha
    # C: Direccion (plain text)
    Set-StringValue $row 3 $r.Direccion
    # D: Comuna (numeric-looking text)
    Set-TextValue $row 4 $r.Comuna
    # E: OT (numeric-looking text)
    Set-TextValue $row 5 $r.OT
    # F: Proveedor Asignado
    Set-StringValue $row 6 $r.Proveedor
    # G: Estado
    Set-StringValue $row 7 $r.Estado
    # H: Observaciones
    Set-StringValue $row 8 $r.Observaciones
    # I: Attachments (number)
    Set-NumberValue $row 9 $r.Attachments
    # J: Tipo de tarea
    Set-StringValue $row 10 $r.Tipo
    # K: Equipo
    Set-StringValue $row 11 $r.Equipo
    # L: Tipo de Elemento
    Set-StringValue $row 12 $r.TipoElemento
    # M: Coordenada_X (number)
    Set-NumberValue $row 13 $r.CoordX
    # N: Coordenada_Y (number)
    Set-NumberValue $row 14 $r.CoordY
    # O: Operacion
    Set-StringValue $row 15 $r.Operacion
    # P: Zona
    Set-StringValue $row 16 $r.Zona
    # Q: PD
    Set-StringValue $row 17 $r.PD
    # R: N2
    Set-StringValue $row 18 $r.N2
}

Write-Host "Added rows 87-88 to sheet $($ws.Name)"
